# ----------------------------------------------------------------------
# Word COM-interop script reproducing the XML diff:
#  1. Append the new "Knarot" (Knärot) section + its reference list after
#     the last paragraph of the document body
#     ("BILAGA 1 - Fridlysta arter").
#  2. Update the date in the first-page header from 2023-09-13 to
#     2023-09-15.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Update the date shown in the first-page header ---
# wdHeaderFooterFirstPage = 2
$hdr = $d.Sections.Item(1).Headers.Item(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

# --- 2. Append the new paragraphs after the last paragraph in the body ---
$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$pendingItalicRanges = @()

# Paragraph 1/13 (style: Heading 1)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Heading 1"
$cur.InsertAfter("Knärot – ekologi samt krav på livsmiljön")

# Paragraph 2/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).")

# Paragraph 3/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Samuel Johnsons doktorsavhandling ")
$rStart = $cur.End - 1
$cur.InsertAfter("“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter(" (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ")
$rStart = $cur.End - 1
$cur.InsertAfter("“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter("Vidare ")
$rStart = $cur.End - 1
$cur.InsertAfter("“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)

# Paragraph 4/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ")
$rStart = $cur.End - 1
$cur.InsertAfter("“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)

# Paragraph 5/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).")

# Paragraph 6/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).")

# Paragraph 7/13 (style: Heading 2)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Heading 2"
$cur.InsertAfter("Referenser - knärot")

# Paragraph 8/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("de Graaf M & Roberts M.R., 2009. ")
$rStart = $cur.End - 1
$cur.InsertAfter("Short-term response of the herbaceous layer within leave patches after harvest. ")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter("Forest Ecology and Management 257, 1014-1025")

# Paragraph 9/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ")
$rStart = $cur.End - 1
$cur.InsertAfter("Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter("Ecological Applications, 22, 2049-2064 ")

# Paragraph 10/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ")
$rStart = $cur.End - 1
$cur.InsertAfter("Interactive effects of drought and edge exposure on old-growth forest understory species. ")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter("Landscape Ecology, 37, sid 1839-1853")

# Paragraph 11/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ")
$rStart = $cur.End - 1
$cur.InsertAfter("Biological legacies buffer local species extinction after logging. ")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter("Journal of Applied Ecology. 51, 53-62.")

# Paragraph 12/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("Skogsstyrelsen, 2022. ")
$rStart = $cur.End - 1
$cur.InsertAfter("Vägledning för hänsyn till knärot. ")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter("https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/")

# Paragraph 13/13 (style: Normal)
$cur.InsertParagraphAfter()

# Stamp italic on any pending ranges from the previous paragraph now
# that the paragraph break has been inserted (keeps the new paragraph
# mark from inheriting italic formatting).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()

$cur = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Paragraphs.Item($d.Paragraphs.Count).Style = "Normal"
$cur.InsertAfter("SLU Artdatabanken, 2021. ")
$rStart = $cur.End - 1
$cur.InsertAfter("Artfaktablad. Naturvård – artfakta. ")
$rEnd = $cur.End - 1
$pendingItalicRanges += ,@($rStart, $rEnd)
$cur.InsertAfter("SLU Artdatabanken, Uppsala ")

# Stamp italic for the very last paragraph (no further paragraph follows).
foreach ($rg in $pendingItalicRanges) {
    $d.Range($rg[0], $rg[1]).Font.Italic = $true
}
$pendingItalicRanges = @()
